# chore: adapt column header formatting to respective input file names
#
# - Rename the "*_old" / "*_new" column header suffixes used throughout the
#   sheet to "*_FV2310" / "*_FV2404" respectively (the AHB format-version
#   identifiers that replaced the generic "old"/"new" labels).
# - Turn the populated range into a proper Excel Table ("Table1") with an
#   AutoFilter, so the header row can be filtered/sorted like a table.
# - Freeze the header row (row 1) so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells -------------------------------------------------
# Columns A1:U1 hold the header labels. Every label ends in either "_old" or
# "_new" (with the sole exception of the "diff" column) - swap those suffixes
# for the new format-version based ones.
$lastCol = 21   # column U
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $headerName = $cell.Value2
    if ($headerName -like "*_old") {
        $cell.Value = ($headerName -replace "_old$", "_FV2310")
    } elseif ($headerName -like "*_new") {
        $cell.Value = ($headerName -replace "_new$", "_FV2404")
    }
}

# --- 2. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Convert the used range into an Excel Table --------------------------
$tableRange = $ws.Range("A1:U82")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

Write-Host "Headers renamed, pane frozen and Table1 created."
